$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.011.67"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "1.907.53"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'318.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "'0.4826"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").Value = "'0.3802"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.07373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'0.9329"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'20.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'0.07751"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "1.923.51"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("D14").Value = "'5.484"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "'6.635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'91.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'0.000008884"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "28.029.67"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'5.137"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "2.140.34"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("D25").Value = "'156.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'1.912"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'2.122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.72%  "
$ws.Range("D29").Value = "'117.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").Value = "'4.981"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").Value = "'3.280"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("D34").Value = "'0.7702"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("D35").Value = "'4.670"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").Value = "'2.588"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.04%  "
$ws.Range("D37").Value = "'0.02054"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "'1.106"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "'0.5497"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "'0.05279"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'3.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "'6.959"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").Value = "'0.1528"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'8.492"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "'110.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.40%  "
$ws.Range("D46").Value = "'10.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "'0.4822"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'1.646"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'67.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "'0.06076"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
